$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 3746283472638430
$ws.Columns.Item(4).ColumnWidth = 11.14
